$wb = $excel.ActiveWorkbook

# --- Add the "metadata" sheet after "data", mirroring the refined scraper output ---
$dataSheet = $wb.Worksheets.Item("data")
$meta = $wb.Worksheets.Add($null, $dataSheet)
$meta.Name = "metadata"

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Aortopathy_Connective Tissue Disorders"
$meta.Range("C2").Value = 44

# "1.55" is written as plain text (not a number) in the source file, so force
# text entry via a temporary "@" number format, then reset the cell format by
# pasting formats from a never-touched cell - this avoids minting a brand new
# (permanent) text-format style index on D2, matching the unstyled target cell.
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.55"
$meta.Range("Z50").Copy()
$meta.Range("D2").PasteSpecial(-4122)

$meta.Range("E2").Value = "2021-09-26T02:48:59.166982Z"
$meta.Range("F2").Value = "2021-10-05 14:33:09.640465"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/44/?format=json"

# Reuse the bold/bordered header style already in the workbook (style index 1)
# for the new sheet's header row and the A2 row-index cell, instead of minting
# a brand-new (duplicate) style via Font/Border property writes.
$dataSheet.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)

$dataSheet.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$meta.Range("A2").Value = 0

# --- Refresh the per-row scrape timestamps on "data" (column F, time_taken) ---
$timestamps = @{
    2 = "2021-10-05 14:33:09.643757"
    3 = "2021-10-05 14:33:09.643765"
    4 = "2021-10-05 14:33:09.643768"
    5 = "2021-10-05 14:33:09.643771"
    6 = "2021-10-05 14:33:09.643774"
    7 = "2021-10-05 14:33:09.643776"
    8 = "2021-10-05 14:33:09.643779"
    9 = "2021-10-05 14:33:09.643781"
    10 = "2021-10-05 14:33:09.643784"
    11 = "2021-10-05 14:33:09.643786"
    12 = "2021-10-05 14:33:09.643789"
    13 = "2021-10-05 14:33:09.643791"
    14 = "2021-10-05 14:33:09.643794"
    15 = "2021-10-05 14:33:09.643796"
    16 = "2021-10-05 14:33:09.643799"
    17 = "2021-10-05 14:33:09.643801"
    18 = "2021-10-05 14:33:09.643804"
    19 = "2021-10-05 14:33:09.643806"
    20 = "2021-10-05 14:33:09.643809"
    21 = "2021-10-05 14:33:09.643811"
    22 = "2021-10-05 14:33:09.643814"
    23 = "2021-10-05 14:33:09.643816"
    24 = "2021-10-05 14:33:09.643818"
    25 = "2021-10-05 14:33:09.643820"
    26 = "2021-10-05 14:33:09.643823"
    27 = "2021-10-05 14:33:09.643826"
    28 = "2021-10-05 14:33:09.643828"
    29 = "2021-10-05 14:33:09.643831"
    30 = "2021-10-05 14:33:09.643833"
    31 = "2021-10-05 14:33:09.643835"
    32 = "2021-10-05 14:33:09.643838"
    33 = "2021-10-05 14:33:09.643840"
    34 = "2021-10-05 14:33:09.643843"
    35 = "2021-10-05 14:33:09.643845"
    36 = "2021-10-05 14:33:09.643848"
    37 = "2021-10-05 14:33:09.643850"
    38 = "2021-10-05 14:33:09.643853"
    39 = "2021-10-05 14:33:09.643855"
    40 = "2021-10-05 14:33:09.643858"
    41 = "2021-10-05 14:33:09.643860"
    42 = "2021-10-05 14:33:09.643863"
    43 = "2021-10-05 14:33:09.643865"
    44 = "2021-10-05 14:33:09.643868"
    45 = "2021-10-05 14:33:09.643871"
    46 = "2021-10-05 14:33:09.643873"
    47 = "2021-10-05 14:33:09.643876"
    48 = "2021-10-05 14:33:09.643878"
    49 = "2021-10-05 14:33:09.643880"
    50 = "2021-10-05 14:33:09.643883"
    51 = "2021-10-05 14:33:09.643885"
    52 = "2021-10-05 14:33:09.643888"
    53 = "2021-10-05 14:33:09.643890"
    54 = "2021-10-05 14:33:09.643893"
    55 = "2021-10-05 14:33:09.643895"
    56 = "2021-10-05 14:33:09.643898"
    57 = "2021-10-05 14:33:09.643900"
    58 = "2021-10-05 14:33:09.643902"
    59 = "2021-10-05 14:33:09.643905"
    60 = "2021-10-05 14:33:09.643907"
    61 = "2021-10-05 14:33:09.643910"
    62 = "2021-10-05 14:33:09.643912"
    63 = "2021-10-05 14:33:09.643915"
    64 = "2021-10-05 14:33:09.643917"
    65 = "2021-10-05 14:33:09.643919"
    66 = "2021-10-05 14:33:09.643923"
    67 = "2021-10-05 14:33:09.643926"
    68 = "2021-10-05 14:33:09.643929"
    69 = "2021-10-05 14:33:09.643931"
    70 = "2021-10-05 14:33:09.643934"
    71 = "2021-10-05 14:33:09.643936"
    72 = "2021-10-05 14:33:09.643939"
    73 = "2021-10-05 14:33:09.643942"
    74 = "2021-10-05 14:33:09.643944"
    75 = "2021-10-05 14:33:09.643947"
    76 = "2021-10-05 14:33:09.643949"
    77 = "2021-10-05 14:33:09.643952"
    78 = "2021-10-05 14:33:09.643957"
    79 = "2021-10-05 14:33:09.643960"
    80 = "2021-10-05 14:33:09.643963"
    81 = "2021-10-05 14:33:09.643965"
    82 = "2021-10-05 14:33:09.643968"
    83 = "2021-10-05 14:33:09.643970"
    84 = "2021-10-05 14:33:09.643973"
    85 = "2021-10-05 14:33:09.643975"
    86 = "2021-10-05 14:33:09.643978"
}
foreach ($row in $timestamps.Keys) {
    $dataSheet.Cells.Item($row, 6).Value = $timestamps[$row]
}

$dataSheet.Activate()
